$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "title"
$ws.Range("B1").Value = "subtitle"
$ws.Range("C1").Value = "bedrooms"
$ws.Range("D1").Value = "price"
$ws.Range("E1").Value = "rating"
$ws.Range("F1").Value = "superhost"
$ws.Range("A2").Value = "Casa em Araruama"
$ws.Range("B2").Value = "Região dos lagos - Araruama - Saquarema - Arraial"
$ws.Range("C2").Value = "2 camas de casal"
$ws.Range("D2").Value = "R$110 por noite"
$ws.Range("E2").Value = "4,79 (282)"
$ws.Range("A3").Value = "Casa em Parque Hotel"
$ws.Range("B3").Value = "Aconchegante Casa - Quintal & Garagem"
$ws.Range("C3").Value = "1 cama queen"
$ws.Range("D3").Value = "R$137 por noite"
$ws.Range("E3").Value = "4,89 (46)"
$ws.Range("F3").Value = "Superhost"
$ws.Range("A4").Value = "Casa em Araruama"
$ws.Range("B4").Value = "Casa em Araruama 1"
$ws.Range("C4").Value = "4 camas"
$ws.Range("D4").Value = "R$106 por noite"
$ws.Range("E4").Value = "4,97 (158)"
$ws.Range("F4").Value = "Superhost"
$ws.Range("A5").Value = "Casa em Araruama"
$ws.Range("B5").Value = "Pedacinho do Céu"
$ws.Range("C5").Value = "1 cama de casal"
$ws.Range("D5").Value = "R$96 por noite"
$ws.Range("E5").Value = "4,92 (12)"
$ws.Range("F5").Value = "Superhost"
$ws.Range("A6").Value = "Casa em Araruama"
$ws.Range("B6").Value = "Recanto para curtir e relaxar em Araruama"
$ws.Range("C6").Value = "3 camas"
$ws.Range("D6").Value = "R$210 por noite"
$ws.Range("E6").Value = "5,0 (29)"
$ws.Range("F6").Value = "Superhost"
$ws.Range("A7").Value = "Apartamento em Araruama"
$ws.Range("B7").Value = "Apartamento em frente a lagoa de Araruama"
$ws.Range("C7").Value = "2 camas"
$ws.Range("D7").Value = "R$173 por noite, originalmente R$211"
$ws.Range("E7").Value = "4,85 (26)"
$ws.Range("F7").ClearContents() | Out-Null
$ws.Range("A8").Value = "Casa em Araruama"
$ws.Range("B8").Value = "Casa tranquila, bem localizada Ar opcional Tv wifi"
$ws.Range("C8").Value = "2 camas"
$ws.Range("D8").Value = "R$118 por noite"
$ws.Range("E8").Value = "4,9 (21)"
$ws.Range("F8").Value = "Superhost"
$ws.Range("A9").Value = "Casa em Araruama"
$ws.Range("B9").Value = "loft <SPA< onde você descansa relaxa e se renova"
$ws.Range("C9").Value = "2 camas"
$ws.Range("D9").Value = "R$189 por noite, originalmente R$302"
$ws.Range("E9").Value = "5,0 (9)"
$ws.Range("F9").Value = "Superhost"
$ws.Range("A10").Value = "Casa em Araruama"
$ws.Range("B10").Value = "Meu aconchego"
$ws.Range("C10").Value = "9 camas"
$ws.Range("D10").Value = "R$187 por noite, originalmente R$235"
$ws.Range("E10").Value = "4,89 (9)"
$ws.Range("F10").ClearContents() | Out-Null
$ws.Range("A11").Value = "Casa em Fazendinha"
$ws.Range("B11").Value = "Casa em Araruama 2"
$ws.Range("C11").Value = "4 camas"
$ws.Range("D11").Value = "R$106 por noite"
$ws.Range("E11").Value = "5,0 (26)"
$ws.Range("F11").Value = "Superhost"
$ws.Range("A12").Value = "Casa de campo em Praia Seca"
$ws.Range("B12").Value = "Casa Maravilhosa com Lagoa privativa"
$ws.Range("C12").Value = "3 camas de casal"
$ws.Range("D12").Value = "R$117 por noite"
$ws.Range("E12").Value = "4,75 (20)"
$ws.Range("A13").Value = "Chalé em Coqueiral"
$ws.Range("B13").Value = "Chalé Recantinho da Lagoa, Araruama,Rj"
$ws.Range("C13").Value = "3 camas"
$ws.Range("D13").Value = "R$91 por noite"
$ws.Range("E13").Value = "4,84 (75)"
$ws.Range("F13").ClearContents() | Out-Null
$ws.Range("A14").Value = "Casa de campo em Outeiro"
$ws.Range("B14").Value = "Casa com piscina no Condomínio Sonho de Vida"
$ws.Range("C14").Value = "4 camas"
$ws.Range("D14").Value = "R$290 por noite, originalmente R$352"
$ws.Range("E14").Value = "5,0 (21)"
$ws.Range("F14").Value = "Superhost"
$ws.Range("A15").Value = "Apartamento em Araruama"
$ws.Range("B15").Value = "Apto na Região dos lagos. Aconchegante e central."
$ws.Range("C15").Value = "2 camas"
$ws.Range("D15").Value = "R$140 por noite"
$ws.Range("E15").Value = "4,71 (7)"
$ws.Range("A16").Value = "Quarto privativo em Araruama"
$ws.Range("B16").Value = "Suíte aconchegante no centro de Praia Seca."
$ws.Range("C16").Value = "1 cama de casal"
$ws.Range("D16").Value = "R$136 por noite, originalmente R$166"
$ws.Range("E16").Value = "5,0 (11)"
$ws.Range("F16").Value = "Superhost"
$ws.Range("A17").Value = "Casa em Parque Hotel"
$ws.Range("B17").Value = "Curta o melhor de Araruama e Região dos Lagos"
$ws.Range("C17").Value = "4 camas"
$ws.Range("D17").Value = "R$133 por noite"
$ws.Range("E17").Value = "4,89 (28)"
$ws.Range("F17").Value = "Superhost"
$ws.Range("A18").Value = "Quarto privativo em Araruama"
$ws.Range("B18").Value = "VEM! Casa Rodrigues - Praia Seca, Araruama"
$ws.Range("C18").Value = "1 cama de casal"
$ws.Range("D18").Value = "R$116 por noite, originalmente R$137"
$ws.Range("E18").Value = "4,92 (173)"
$ws.Range("F18").Value = "Superhost"
$ws.Range("A19").Value = "Casa em Araruama"
$ws.Range("B19").Value = "Aluguel por temporada com piscina privativa"
$ws.Range("C19").Value = "3 camas"
$ws.Range("D19").Value = "R$315 por noite"
$ws.Range("E19").Value = "4,97 (66)"
$ws.Range("F19").Value = "Superhost"
